$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''67.691.08'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = '''3.503.91'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").Value = '''605.69'
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("D6").Value = '''150.28'
$ws.Range("E6").Value = '  +0.47%  '
$ws.Range("D7").Value = '''3.506.13'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '''0.484'
$ws.Range("E9").Value = '  +0.69%  '
$ws.Range("D10").Value = '''0.142'
$ws.Range("E10").Value = '  +2.68%  '
$ws.Range("D11").Value = '''7.48'
$ws.Range("E11").Value = '  +6.40%  '
$ws.Range("D12").Value = '''0.428'
$ws.Range("E12").Value = '  +1.08%  '
$ws.Range("D13").Value = '''32.31'
$ws.Range("E13").Value = '  +2.50%  '
$ws.Range("E14").Value = '  -1.79%  '
$ws.Range("D15").Value = '''4.099.81'
$ws.Range("E15").Value = '  +0.18%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '''68.107.61'
$ws.Range("E16").Value = '  +1.02%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '''3.505.24'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").Value = '''6.50'
$ws.Range("E19").Value = '  +1.33%  '
$ws.Range("D20").Value = '''15.40'
$ws.Range("E20").Value = '  +1.67%  '
$ws.Range("D21").Value = '''9.67'
$ws.Range("E21").Value = '  +6.03%  '
$ws.Range("D22").Value = '''446.18'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = '''0.626'
$ws.Range("E23").Value = '  +0.65%  '
$ws.Range("D24").Value = '''77.72'
$ws.Range("E24").Value = '  +0.47%  '
$ws.Range("D25").Value = '''3.649.22'
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  +0.48%  '
$ws.Range("D28").Value = '''8.78'
$ws.Range("E28").Value = '  +5.77%  '
$ws.Range("D29").Value = '''9.95'
$ws.Range("E29").Value = '  -3.61%  '
$ws.Range("E30").Value = '  +0.47%  '
$ws.Range("D31").Value = '''1.63'
$ws.Range("E31").Value = '  +6.25%  '
$ws.Range("D32").Value = '''0.167'
$ws.Range("E32").Value = '  +1.63%  '
$ws.Range("E33").Value = '  +0.19%  '
$ws.Range("D34").Value = '''25.52'
$ws.Range("E34").Value = '  -0.57%  '
$ws.Range("D35").Value = '''6.12'
$ws.Range("E35").Value = '  +0.62%  '
$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D36").Value = '''3.499.61'
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''1.85'
$ws.Range("E37").Value = '  +1.00%  '
$ws.Range("D38").Value = '''7.93'
$ws.Range("E38").Value = '  -1.61%  '
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("D40").Value = '''2.29'
$ws.Range("E40").Value = '  +5.30%  '
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("D42").Value = '''175.54'
$ws.Range("E42").Value = '  -0.99%  '
$ws.Range("D43").Value = '''0.0887'
$ws.Range("E43").Value = '  +1.76%  '
$ws.Range("D44").Value = '''5.43'
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '''0.879'
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '''29.82'
$ws.Range("E46").Value = '  +8.54%  '
$ws.Range("D47").Value = '''46.56'
$ws.Range("E47").Value = '  +2.70%  '
$ws.Range("D48").Value = '''1.29'
$ws.Range("E48").Value = '  +3.52%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = '''7.59'
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").Value = '''2.50'
$ws.Range("E50").Value = '  -1.70%  '
$ws.Range("D51").Value = '''0.251'
$ws.Range("E51").Value = '  +2.68%  '
